$d = $word.ActiveDocument

# 1. "2487097" -> "2445987" (2 occurrences, identical text/context)
$r = $d.Content
$r.Find.Execute("2487097", $true, $true, $false, $false, $false, $true, 1, $false, "2445987", 2)

# 2. "03.05.2024" -> "01.05.2024" (2 occurrences)
$r = $d.Content
$r.Find.Execute("03.05.2024", $true, $true, $false, $false, $false, $true, 1, $false, "01.05.2024", 2)

# 3. "Светлый" -> "Светлый " (2 occurrences, add trailing space)
$r = $d.Content
$r.Find.Execute("Светлый", $true, $true, $false, $false, $false, $true, 1, $false, "Светлый ", 2)

# 4. "Кодекса торгового мореплавания (КТМ РФ)" -> add trailing space (1 occurrence)
$r = $d.Content
$r.Find.Execute("Кодекса торгового мореплавания (КТМ РФ)", $true, $true, $false, $false, $false, $true, 1, $false, "Кодекса торгового мореплавания (КТМ РФ) ", 2)

# 5. "ВОЛГА" -> "СИНЕГОРСК"
$r = $d.Content
$r.Find.Execute("ВОЛГА", $true, $true, $false, $false, $false, $true, 1, $false, "СИНЕГОРСК", 2)

# 6. "940330" -> "021026"
$r = $d.Content
$r.Find.Execute("940330", $true, $true, $false, $false, $false, $true, 1, $false, "021026", 2)

# 7. " освидетельствование" -> "Очередное освидетельствование "
$r = $d.Content
$r.Find.Execute(" освидетельствование", $true, $true, $false, $false, $false, $true, 1, $false, "Очередное освидетельствование ", 2)

# 8. "Свидетельство ф. 6.5.30 №№ 24.42.03.00234.121 - 24.42.03.00236.121 от --" -> "Согласно перечню ф. 6.4.7-1 № 121-212-08-343489 от --"
$r = $d.Content
$r.Find.Execute("Свидетельство ф. 6.5.30 №№ 24.42.03.00234.121 - 24.42.03.00236.121 от --", $true, $true, $false, $false, $false, $true, 1, $false, "Согласно перечню ф. 6.4.7-1 № 121-212-08-343489 от --", 2)

# 9. "10 000,00 p. (десять тысяч рублей 00 копеек)" -> "5 150,00 p. (пять тысяч сто пятьдесят рублей 00 копеек)"
$r = $d.Content
$r.Find.Execute("10 000,00 p. (десять тысяч рублей 00 копеек)", $true, $true, $false, $false, $false, $true, 1, $false, "5 150,00 p. (пять тысяч сто пятьдесят рублей 00 копеек)", 2)
